# Append the two new daily GSC export rows (2025-11-30, 2025-12-01) to the
# "Chart" worksheet, matching the pattern of the existing rows (0 non-HTTPS
# URLs, 27 HTTPS URLs / pages).
#
# NOTE: writing an ISO date string ("2025-11-30") straight into .Value gets
# auto-recognised as a real date by the engine and re-stamped with a brand
# new date-formatted style, which would also touch styles.xml (not part of
# this change). Routing the text through a TEXT() formula + paste-as-values
# keeps the literal string (shared-string, default style) like the other
# Date column cells.
$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

function Set-TextValue($cell, $text) {
    $cell.Formula = '=TEXT("' + $text + '","yyyy-mm-dd")'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

Set-TextValue $chart.Cells.Item(56, 1) "2025-11-30"
$chart.Cells.Item(56, 2).Value = 0.0
$chart.Cells.Item(56, 3).Value = 27.0

Set-TextValue $chart.Cells.Item(57, 1) "2025-12-01"
$chart.Cells.Item(57, 2).Value = 0.0
$chart.Cells.Item(57, 3).Value = 27.0
